$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to C4
$ws.Range("C4").Select()

# Update the cost value in C3
$ws.Range("C3").Value = 19.989999999999998

# Update the hyperlink formula in D3
$ws.Range("D3").Formula = '=HYPERLINK("https://www.amazon.com/Inland-1-75mm-Natural-Printer-Filament/dp/B00YSP5SR8/ref=sr_1_3?dchild=1&keywords=pla+filament+natural&qid=1618415287&sr=8-3","Filament")'
